$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B133").Value = 7515550
$ws.Range("E133").Value = "Gil Vicente"
$ws.Range("F133").Value = "Boavista"
$ws.Range("G133").Value = 1
$ws.Range("I133").Value = 0
$ws.Range("L133").Value = 2.3
$ws.Range("M133").Value = 3.5
$ws.Range("N133").Value = 2.9
$ws.Range("O133").Value = 2.3
$ws.Range("P133").Value = 3.3
$ws.Range("Q133").Value = 3
$ws.Range("R133").Value = -0.25
$ws.Range("S133").Value = 2.05
$ws.Range("T133").Value = 1.8
$ws.Range("U133").Value = 2.5
$ws.Range("V133").Value = 2.05
$ws.Range("W133").Value = 1.8
$ws.Range("X133").Value = 1.3
$ws.Range("AA133").Value = 1.05
$ws.Range("AC133").Value = -1
$ws.Range("AD133").Value = 0.8
$ws.Range("B134").Value = 7513577
$ws.Range("E134").Value = "Estoril"
$ws.Range("F134").Value = "SC Farense"
$ws.Range("G134").Value = 4
$ws.Range("I134").Value = 2
$ws.Range("L134").Value = 2.15
$ws.Range("M134").Value = 3.6
$ws.Range("N134").Value = 3.2
$ws.Range("O134").Value = 1.833
$ws.Range("P134").Value = 4
$ws.Range("Q134").Value = 3.8
$ws.Range("R134").Value = -0.5
$ws.Range("S134").Value = 1.875
$ws.Range("T134").Value = 1.975
$ws.Range("U134").Value = 2.75
$ws.Range("V134").Value = 1.875
$ws.Range("W134").Value = 1.975
$ws.Range("X134").Value = 0.833
$ws.Range("AA134").Value = 0.875
$ws.Range("AC134").Value = 0.875
$ws.Range("AD134").Value = -1
$ws.Range("B151").Value = 6876591
$ws.Range("E151").Value = "Vizela"
$ws.Range("F151").Value = "Boavista"
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 4
$ws.Range("I151").Value = 0
$ws.Range("J151").Value = 2
$ws.Range("K151").Value = "A"
$ws.Range("L151").Value = 2.3
$ws.Range("M151").Value = 3.1
$ws.Range("N151").Value = 3.25
$ws.Range("O151").Value = 1.95
$ws.Range("P151").Value = 3.2
$ws.Range("Q151").Value = 4.2
$ws.Range("R151").Value = -0.5
$ws.Range("U151").Value = 2.25
$ws.Range("V151").Value = 2
$ws.Range("W151").Value = 1.85
$ws.Range("X151").Value = -1
$ws.Range("Z151").Value = 3.2
$ws.Range("AA151").Value = -1
$ws.Range("AB151").Value = 0.825
$ws.Range("AC151").Value = 1
$ws.Range("B152").Value = 6876586
$ws.Range("E152").Value = "Benfica"
$ws.Range("F152").Value = "Rio Ave"
$ws.Range("G152").Value = 4
$ws.Range("H152").Value = 1
$ws.Range("I152").Value = 1
$ws.Range("J152").Value = 1
$ws.Range("K152").Value = "H"
$ws.Range("L152").Value = 1.166
$ws.Range("M152").Value = 7.5
$ws.Range("N152").Value = 15
$ws.Range("O152").Value = 1.2
$ws.Range("P152").Value = 8
$ws.Range("Q152").Value = 10
$ws.Range("R152").Value = -2
$ws.Range("U152").Value = 3.25
$ws.Range("V152").Value = 1.925
$ws.Range("W152").Value = 1.925
$ws.Range("X152").Value = 0.2
$ws.Range("Z152").Value = -1
$ws.Range("AA152").Value = 1.025
$ws.Range("AB152").Value = -1
$ws.Range("AC152").Value = 0.925
$ws.Range("B195").Value = 6876630
$ws.Range("E195").Value = "Benfica"
$ws.Range("F195").Value = "Vizela"
$ws.Range("G195").Value = 6
$ws.Range("H195").Value = 1
$ws.Range("I195").Value = 5
$ws.Range("K195").Value = "H"
$ws.Range("L195").Value = 1.111
$ws.Range("M195").Value = 8.5
$ws.Range("N195").Value = 21
$ws.Range("O195").Value = 1.1
$ws.Range("P195").Value = 9.5
$ws.Range("Q195").Value = 23
$ws.Range("R195").Value = -2
$ws.Range("S195").Value = 1.89
$ws.Range("T195").Value = 2.01
$ws.Range("U195").Value = 3.25
$ws.Range("V195").Value = 1.825
$ws.Range("W195").Value = 2.025
$ws.Range("X195").Value = 0.1000000000000001
$ws.Range("Z195").Value = -1
$ws.Range("AA195").Value = 0.8899999999999999
$ws.Range("AB195").Value = -1
$ws.Range("AC195").Value = 0.825
$ws.Range("B196").Value = 6876633
$ws.Range("E196").Value = "Estoril"
$ws.Range("F196").Value = "Gil Vicente"
$ws.Range("G196").Value = 1
$ws.Range("H196").Value = 3
$ws.Range("I196").Value = 0
$ws.Range("K196").Value = "A"
$ws.Range("L196").Value = 2.2
$ws.Range("M196").Value = 3.4
$ws.Range("N196").Value = 3.2
$ws.Range("O196").Value = 2.1
$ws.Range("P196").Value = 3.5
$ws.Range("Q196").Value = 3.5
$ws.Range("R196").Value = -0.25
$ws.Range("S196").Value = 1.95
$ws.Range("T196").Value = 1.95
$ws.Range("U196").Value = 2.5
$ws.Range("V196").Value = 1.925
$ws.Range("W196").Value = 1.925
$ws.Range("X196").Value = -1
$ws.Range("Z196").Value = 2.5
$ws.Range("AA196").Value = -1
$ws.Range("AB196").Value = 0.95
$ws.Range("AC196").Value = 0.925

Write-Host "Applied Portugal Primeira Liga updates"
